$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.542.97"
$ws.Range("E2").Value = "  +0.40%  "

# Row 3
$ws.Range("D3").Value = "3.019.99"
$ws.Range("E3").Value = "  +0.38%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "'595.45"
$ws.Range("E5").Value = "  +2.14%  "

# Row 6
$ws.Range("D6").Value = "'147.86"
$ws.Range("E6").Value = "  +1.35%  "

# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").Value = "3.018.31"
$ws.Range("E8").Value = "  +0.38%  "

# Row 9
$ws.Range("E9").Value = "  -2.03%  "

# Row 10
$ws.Range("D10").Value = "'6.28"
$ws.Range("E10").Value = "  +8.43%  "

# Row 11
$ws.Range("E11").Value = "  +0.88%  "

# Row 12
$ws.Range("E12").Value = "  -2.13%  "

# Row 13
$ws.Range("D13").Value = "'0.0000233"
$ws.Range("E13").Value = "  +1.87%  "

# Row 14
$ws.Range("D14").Value = "'34.52"
$ws.Range("E14").Value = "  +0.10%  "

# Row 15
$ws.Range("E15").Value = "  +2.61%  "

# Row 16
$ws.Range("D16").Value = "3.520.13"
$ws.Range("E16").Value = "  +0.51%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "62.477.95"
$ws.Range("E17").Value = "  +0.30%  "

# Row 18
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.01"
$ws.Range("E18").Value = "  -1.92%  "

# Row 19
$ws.Range("D19").Value = "3.019.92"
$ws.Range("E19").Value = "  +0.40%  "

# Row 20
$ws.Range("D20").Value = "'450.31"
$ws.Range("E20").Value = "  -2.01%  "

# Row 21
$ws.Range("D21").Value = "'14.13"
$ws.Range("E21").Value = "  +1.06%  "

# Row 22
$ws.Range("D22").Value = "'0.687"
$ws.Range("E22").Value = "  -0.16%  "

# Row 23
$ws.Range("D23").Value = "'7.38"
$ws.Range("E23").Value = "  -0.44%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'82.01"
$ws.Range("E24").Value = "  +0.20%  "

# Row 25
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "'11.18"
$ws.Range("E25").Value = "  +11.62%  "

# Row 26
$ws.Range("D26").Value = "'2.28"
$ws.Range("E26").Value = "  +3.43%  "

# Row 27
$ws.Range("E27").Value = "  -0.48%  "

# Row 28
$ws.Range("E28").Value = "  +0.00%  "

# Row 29
$ws.Range("D29").Value = "'2.71"
$ws.Range("E29").Value = "  +3.88%  "

# Row 30
$ws.Range("D30").Value = "'7.33"
$ws.Range("E30").Value = "  +5.21%  "

# Row 31
$ws.Range("E31").Value = "  -0.20%  "

# Row 32
$ws.Range("E32").Value = "  +1.12%  "

# Row 33
$ws.Range("D33").Value = "'27.50"
$ws.Range("E33").Value = "  -3.55%  "

# Row 34
$ws.Range("E34").Value = "  +2.08%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0853"
$ws.Range("E35").Value = "  +7.03%  "

# Row 36
$ws.Range("D36").Value = "'1.02"
$ws.Range("E36").Value = "  -0.42%  "

# Row 37
$ws.Range("D37").Value = "'5.84"
$ws.Range("E37").Value = "  +1.15%  "

# Row 38
$ws.Range("D38").Value = "'50.45"
$ws.Range("E38").Value = "  +0.30%  "

# Row 39
$ws.Range("E39").Value = "  -1.66%  "

# Row 40
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'2.98"
$ws.Range("E40").Value = "  +3.45%  "

# Row 41
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'8.99"
$ws.Range("E41").Value = "  -2.35%  "

# Row 42
$ws.Range("D42").Value = "'0.125"
$ws.Range("E42").Value = "  +6.95%  "

# Row 43
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'406.94"
$ws.Range("E43").Value = "  +3.75%  "

# Row 44
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").Value = "'41.28"
$ws.Range("E44").Value = "  +12.63%  "

# Row 45
$ws.Range("D45").Value = "'0.280"
$ws.Range("E45").Value = "  +4.29%  "

# Row 46
$ws.Range("D46").Value = "'0.0354"
$ws.Range("E46").Value = "  -0.93%  "

# Row 47
$ws.Range("D47").Value = "2.718.35"
$ws.Range("E47").Value = "  -0.32%  "

# Row 48
$ws.Range("D48").Value = "'132.93"
$ws.Range("E48").Value = "  +3.52%  "

# Row 49
$ws.Range("E49").Value = "  +0.08%  "

# Row 50
$ws.Range("E50").Value = "  -0.03%  "

# Row 51
$ws.Range("E51").Value = "  -1.46%  "
